$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Demi (sheet2): fix the mis-typed date in B18 (was literal text "27/11/20203")
# to a real date, and leave a note in B8 explaining the printing glitch that
# caused the typo in the first place.
# This must run BEFORE the other new-string edits so the freed shared-string
# slot (old index 43) is reused first, matching natural Excel behaviour.
# ---------------------------------------------------------------------------
$demi = $wb.Worksheets.Item("Demi")
$demi.Range("B18").Value = 45257
$demi.Range("B8").Value = "Datum word verkeerd geprint. Staat wel correct in edit balk boven in excel"

# ---------------------------------------------------------------------------
# Marvin (sheet1): remark on the existing last row, plus two new rows for
# yesterday's and today's KBS lessons.
# ---------------------------------------------------------------------------
$marvin = $wb.Worksheets.Item("Marvin")
$marvin.Range("D23").Value = "zoinks scoob"

$marvin.Range("B23").Copy()
$marvin.Range("B24:B25").PasteSpecial(-4122)
$marvin.Range("B24").Value = 45265
$marvin.Range("B25").Value = 45265

$marvin.Cells.Item(24, 1).Value = "KBS b les"
$marvin.Cells.Item(24, 3).Value = 120
$marvin.Cells.Item(24, 4).Value = "Les"

$marvin.Cells.Item(25, 1).Value = "KBS a les"
$marvin.Cells.Item(25, 3).Value = 120
$marvin.Cells.Item(25, 4).Value = "Les"

$marvin.Range("A24:D25").Select()

# ---------------------------------------------------------------------------
# Demi (sheet2): the same two new rows for yesterday's and today's KBS lessons
# ---------------------------------------------------------------------------
$demi.Range("B17").Copy()
$demi.Range("B23:B24").PasteSpecial(-4122)
$demi.Range("B23").Value = 45265
$demi.Range("B24").Value = 45265

$demi.Cells.Item(23, 1).Value = "KBS b les"
$demi.Cells.Item(23, 3).Value = 120
$demi.Cells.Item(23, 4).Value = "Les"

$demi.Cells.Item(24, 1).Value = "KBS a les"
$demi.Cells.Item(24, 3).Value = 120
$demi.Cells.Item(24, 4).Value = "Les"

$demi.Range("A23:D24").Select()

# ---------------------------------------------------------------------------
# Lucas (sheet3): the same two new rows
# ---------------------------------------------------------------------------
$lucas = $wb.Worksheets.Item("Lucas")
$lucas.Range("B21").Copy()
$lucas.Range("B22:B23").PasteSpecial(-4122)
$lucas.Range("B22").Value = 45265
$lucas.Range("B23").Value = 45265

$lucas.Cells.Item(22, 1).Value = "KBS b les"
$lucas.Cells.Item(22, 3).Value = 120
$lucas.Cells.Item(22, 4).Value = "Les"

$lucas.Cells.Item(23, 1).Value = "KBS a les"
$lucas.Cells.Item(23, 3).Value = 120
$lucas.Cells.Item(23, 4).Value = "Les"

$lucas.Range("A22:D23").Select()

# ---------------------------------------------------------------------------
# Luuk (sheet4): the same two new rows, plus a third row noting Jochem
# ("Lonely") joined today's planned KBS moment over Teams.
# ---------------------------------------------------------------------------
$luuk = $wb.Worksheets.Item("Luuk")
$luuk.Range("B24").Copy()
$luuk.Range("B25:B27").PasteSpecial(-4122)
$luuk.Range("B25").Value = 45265
$luuk.Range("B26").Value = 45265
$luuk.Range("B27").Value = 45266

$luuk.Cells.Item(25, 1).Value = "KBS b les"
$luuk.Cells.Item(25, 3).Value = 120
$luuk.Cells.Item(25, 4).Value = "Les"

$luuk.Cells.Item(26, 1).Value = "KBS a les"
$luuk.Cells.Item(26, 3).Value = 120
$luuk.Cells.Item(26, 4).Value = "Les"

$luuk.Cells.Item(27, 1).Value = "KBS momentje geplanned"
$luuk.Cells.Item(27, 3).Value = 120
$luuk.Cells.Item(27, 4).Value = "Lonely aanwezig in teams D:"

# ---------------------------------------------------------------------------
# Jochem (sheet5): the same two new rows
# ---------------------------------------------------------------------------
$jochem = $wb.Worksheets.Item("Jochem")
$jochem.Range("B20").Copy()
$jochem.Range("B21:B22").PasteSpecial(-4122)
$jochem.Range("B21").Value = 45265
$jochem.Range("B22").Value = 45265

$jochem.Cells.Item(21, 1).Value = "KBS b les"
$jochem.Cells.Item(21, 3).Value = 120
$jochem.Cells.Item(21, 4).Value = "Les"

$jochem.Cells.Item(22, 1).Value = "KBS a les"
$jochem.Cells.Item(22, 3).Value = 120
$jochem.Cells.Item(22, 4).Value = "Les"

$jochem.Range("A21:D22").Select()

# ---------------------------------------------------------------------------
# End on Luuk's sheet, with D28 selected (next free cell after the new rows)
# ---------------------------------------------------------------------------
$luuk.Activate()
$luuk.Range("D28").Select()
